$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2073439602192391
$ws.Range("C2").Value = 2.184041328014591
$ws.Range("D2").Value = 14.59467764566177
$ws.Range("E2").Value = 3.820298109527811
$ws.Range("F2").Value = 3.867285822301346
$ws.Range("G2").Value = 37
$ws.Range("B3").Value = 0.6636414647983299
$ws.Range("C3").Value = 1.63741016751905
$ws.Range("D3").Value = 11.25093469053298
$ws.Range("E3").Value = 3.35424129879366
$ws.Range("F3").Value = 3.334574417404926
$ws.Range("G3").Value = 36
$ws.Range("B4").Value = 0.4237612586259278
$ws.Range("C4").Value = 1.679542155191349
$ws.Range("D4").Value = 11.07338574736448
$ws.Range("E4").Value = 3.327669717289334
$ws.Range("F4").Value = 3.3487637095731
$ws.Range("G4").Value = 35
$ws.Range("B5").Value = 0.6888881351581906
$ws.Range("C5").Value = 1.666778803909363
$ws.Range("D5").Value = 11.55801848543534
$ws.Range("E5").Value = 3.399708588310966
$ws.Range("F5").Value = 3.379247518761783
$ws.Range("G5").Value = 34
$ws.Range("B6").Value = 0.4816360917539647
$ws.Range("C6").Value = 1.844218230436959
$ws.Range("D6").Value = 12.08105673389089
$ws.Range("E6").Value = 3.475781456577915
$ws.Range("F6").Value = 3.495621155895251
$ws.Range("G6").Value = 33
$ws.Range("B7").Value = 0.7221805366241528
$ws.Range("C7").Value = 1.78756334845686
$ws.Range("D7").Value = 12.36122690756038
$ws.Range("E7").Value = 3.515853652750692
$ws.Range("F7").Value = 3.495941562969434
$ws.Range("G7").Value = 32
$ws.Range("B8").Value = 0.5718078270650869
$ws.Range("C8").Value = 1.854104342957878
$ws.Range("D8").Value = 12.66500268672434
$ws.Range("E8").Value = 3.558792307331847
$ws.Range("F8").Value = 3.570617095893344
$ws.Range("G8").Value = 31
$ws.Range("B9").Value = 0.7779467227852803
$ws.Range("C9").Value = 1.816601862872559
$ws.Range("D9").Value = 13.04591585095101
$ws.Range("E9").Value = 3.611913045873476
$ws.Range("F9").Value = 3.587437095081128
$ws.Range("G9").Value = 30
$ws.Range("B10").Value = 0.5192813677703352
$ws.Range("C10").Value = 1.862221145406507
$ws.Range("D10").Value = 13.30414656006889
$ws.Range("E10").Value = 3.647484963652199
$ws.Range("F10").Value = 3.674236116982628
$ws.Range("G10").Value = 29
$ws.Range("B11").Value = 0.8142683941004476
$ws.Range("C11").Value = 1.860469755909099
$ws.Range("D11").Value = 13.97793062376009
$ws.Range("E11").Value = 3.738707079159865
$ws.Range("F11").Value = 3.715917378240766
$ws.Range("G11").Value = 28
